# The underlying game-log now includes additional simulated games, which
# changes the empirically-derived state transition probabilities on the
# Starting_State matrix (Sheet1). Re-write the recomputed probabilities
# for every affected "from state" row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Af0 -- updated transition probabilities
$ws.Range("B2").Value = 0.1818181818181818
$ws.Range("C2").Value = 0.5909090909090909
$ws.Range("J2").Value = 0.04545454545454546
$ws.Range("P2").Value = 0.1363636363636364
$ws.Range("S2").Value = 0.04545454545454546

# Row 3: Af1 -- updated transition probabilities
$ws.Range("J3").Value = 0.07692307692307693
$ws.Range("P3").Value = 0.8461538461538461
$ws.Range("S3").Value = 0.07692307692307693

# Row 4: Af2 -- updated transition probabilities
$ws.Range("J4").Value = 0.125
$ws.Range("P4").Value = 0.625
$ws.Range("S4").Value = 0.25

# Row 6: Ai0 -- updated transition probabilities
$ws.Range("B6").Value = 0.07692307692307693
$ws.Range("J6").Value = 0.2307692307692308
$ws.Range("Q6").Value = 0.2307692307692308
$ws.Range("R6").Value = 0.07692307692307693
$ws.Range("S6").Value = 0.3846153846153846

# Row 7: Ai1 -- updated transition probabilities
$ws.Range("B7").Value = 0.09090909090909091
$ws.Range("J7").Value = 0.3636363636363636
$ws.Range("Q7").Value = 0.2727272727272727
$ws.Range("S7").Value = 0.2727272727272727

# Row 8: Ai2 -- updated transition probabilities
$ws.Range("B8").Value = 0.08571428571428572
$ws.Range("D8").Value = 0.05714285714285714
$ws.Range("F8").Value = 0.08571428571428572
$ws.Range("J8").Value = 0.08571428571428572
$ws.Range("O8").Value = 0.02857142857142857
$ws.Range("Q8").Value = 0.2857142857142857
$ws.Range("R8").Value = 0.1142857142857143
$ws.Range("S8").Value = 0.2571428571428571

# Row 9: Ai3 -- updated transition probabilities
$ws.Range("B9").Value = 0.1176470588235294
$ws.Range("F9").Value = 0.05882352941176471
$ws.Range("J9").Value = 0.1176470588235294
$ws.Range("Q9").Value = 0.1764705882352941
$ws.Range("R9").Value = 0.1764705882352941
$ws.Range("S9").Value = 0.3529411764705883

# Row 10: Ar0 -- updated transition probabilities
$ws.Range("B10").Value = 0.08029197080291971
$ws.Range("D10").Value = 0.04379562043795621
$ws.Range("F10").Value = 0.0364963503649635
$ws.Range("J10").Value = 0.0948905109489051
$ws.Range("O10").Value = 0.0145985401459854
$ws.Range("Q10").Value = 0.3284671532846715
$ws.Range("R10").Value = 0.1313868613138686
$ws.Range("S10").Value = 0.2700729927007299

# Row 11: Bf0 -- updated transition probabilities
$ws.Range("G11").Value = 0.2941176470588235
$ws.Range("J11").Value = 0.1176470588235294
$ws.Range("K11").Value = 0.2352941176470588
$ws.Range("L11").Value = 0.3529411764705883

# Row 12: Bf1 -- updated transition probabilities
$ws.Range("G12").Value = 0.6666666666666666
$ws.Range("K12").Value = 0.1666666666666667
$ws.Range("S12").Value = 0.1666666666666667

# Row 15: Bi0 -- updated transition probabilities
$ws.Range("H15").Value = 0.1538461538461539
$ws.Range("I15").Value = 0.07692307692307693
$ws.Range("J15").Value = 0.3846153846153846
$ws.Range("S15").Value = 0.3846153846153846

# Row 16: Bi1 -- updated transition probabilities
$ws.Range("H16").Value = 0.2222222222222222
$ws.Range("I16").Value = 0.2222222222222222
$ws.Range("J16").Value = 0.2777777777777778
$ws.Range("O16").Value = 0.1111111111111111
$ws.Range("S16").Value = 0.1666666666666667

# Row 17: Bi2 -- updated transition probabilities
$ws.Range("H17").Value = 0.09677419354838709
$ws.Range("I17").Value = 0.09677419354838709
$ws.Range("J17").Value = 0.5645161290322581
$ws.Range("K17").Value = 0.04838709677419355
$ws.Range("M17").Value = 0.01612903225806452
$ws.Range("O17").Value = 0.01612903225806452
$ws.Range("S17").Value = 0.1612903225806452

# Row 18: Bi3 -- updated transition probabilities
$ws.Range("H18").Value = 0.1851851851851852
$ws.Range("I18").Value = 0.03703703703703703
$ws.Range("J18").Value = 0.4814814814814815
$ws.Range("O18").Value = 0.1111111111111111
$ws.Range("S18").Value = 0.1851851851851852

# Row 19: Br0 -- updated transition probabilities
$ws.Range("F19").Value = 0.01075268817204301
$ws.Range("H19").Value = 0.1935483870967742
$ws.Range("I19").Value = 0.06451612903225806
$ws.Range("J19").Value = 0.5268817204301075
$ws.Range("K19").Value = 0.09677419354838709
$ws.Range("M19").Value = 0.01075268817204301
$ws.Range("O19").Value = 0.03225806451612903
$ws.Range("S19").Value = 0.06451612903225806

